# Auto-generated edit script: update cryptos price/volume columns (and two swapped rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.850.46'
$ws.Range('E2').Value = '  +0.24%  '
$ws.Range('D3').Value = '1.649.12'
$ws.Range('E3').Value = '  +2.24%  '
$ws.Range('D4').Value = '''1.004'
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').Value = '''308.79'
$ws.Range('E5').Value = '  +0.68%  '
$ws.Range('D7').Value = '''0.3881'
$ws.Range('E7').Value = '  -0.27%  '
$ws.Range('D8').Value = '''0.3832'
$ws.Range('E8').Value = '  +0.91%  '
$ws.Range('D9').Value = '''51.07'
$ws.Range('E9').Value = '  +5.01%  '
$ws.Range('D10').Value = '''1.348'
$ws.Range('E10').Value = '  -0.04%  '
$ws.Range('D11').Value = '''1.004'
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('D12').Value = '''0.08443'
$ws.Range('E12').Value = '  +0.37%  '
$ws.Range('D13').Value = '''23.85'
$ws.Range('E13').Value = '  +0.38%  '
$ws.Range('D14').Value = '''7.129'
$ws.Range('E14').Value = '  +1.70%  '
$ws.Range('D15').Value = '''7.818'
$ws.Range('E15').Value = '  +5.05%  '
$ws.Range('D16').Value = '''0.00001307'
$ws.Range('E16').Value = '  +2.95%  '
$ws.Range('D17').Value = '1.648.62'
$ws.Range('E17').Value = '  +3.15%  '
$ws.Range('D18').Value = '''94.54'
$ws.Range('E18').Value = '  +1.71%  '
$ws.Range('E19').Value = '  +1.31%  '
$ws.Range('D20').Value = '''19.72'
$ws.Range('E20').Value = '  -1.14%  '
$ws.Range('D21').Value = '''6.880'
$ws.Range('E21').Value = '  +1.49%  '
$ws.Range('D23').Value = '''13.56'
$ws.Range('E23').Value = '  +1.31%  '
$ws.Range('D24').Value = '23.850.66'
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('D25').Value = '''2.469'
$ws.Range('E25').Value = '  +1.20%  '
$ws.Range('D26').Value = '''3.035'
$ws.Range('E26').Value = '  +8.48%  '
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').Value = '''152.36'
$ws.Range('E28').Value = '  -2.91%  '
$ws.Range('D29').Value = '''5.447'
$ws.Range('E29').Value = '  +3.74%  '
$ws.Range('D30').Value = '''139.23'
$ws.Range('E30').Value = '  +0.12%  '
$ws.Range('D31').Value = '''7.767'
$ws.Range('E31').Value = '  +0.52%  '
$ws.Range('D32').Value = '''2.506'
$ws.Range('E32').Value = '  +0.82%  '
$ws.Range('D33').Value = '1.834.32'
$ws.Range('E33').Value = '  +2.43%  '
$ws.Range('D34').Value = '''1.022'
$ws.Range('E34').Value = '  +6.09%  '
$ws.Range('D35').Value = '''0.08006'
$ws.Range('E35').Value = '  -1.23%  '
$ws.Range('D36').Value = '''0.02949'
$ws.Range('E36').Value = '  +3.09%  '
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').Value = '''6.653'
$ws.Range('E37').Value = '  +2.11%  '
$ws.Range('B38').Value = 'FraxShare'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D38').Value = '''10.89'
$ws.Range('E38').Value = '  +5.22%  '
$ws.Range('D39').Value = '''0.2678'
$ws.Range('E39').Value = '  +1.41%  '
$ws.Range('D40').Value = '''0.09113'
$ws.Range('E40').Value = '  -0.17%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').Value = '''13.49'
$ws.Range('E41').Value = '  +0.67%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = '''0.7513'
$ws.Range('E42').Value = '  +1.43%  '
$ws.Range('D43').Value = '''1.416'
$ws.Range('E43').Value = '  -0.34%  '
$ws.Range('D44').Value = '''16.27'
$ws.Range('E44').Value = '  +3.69%  '
$ws.Range('D45').Value = '''0.6913'
$ws.Range('E45').Value = '  +1.64%  '
$ws.Range('D46').Value = '''2.450'
$ws.Range('E46').Value = '  +0.85%  '
$ws.Range('D47').Value = '''4.066'
$ws.Range('E47').Value = '  +0.40%  '
$ws.Range('D48').Value = '''1.000'
$ws.Range('E48').Value = '  -0.13%  '
$ws.Range('D49').Value = '''0.08259'
$ws.Range('E49').Value = '  +0.63%  '
$ws.Range('D50').Value = '''133.78'
$ws.Range('E50').Value = '  +1.23%  '
$ws.Range('D51').Value = '''1.224'
$ws.Range('E51').Value = '  +2.38%  '
